$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D6").Value = 436
